$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("B1").Value = "Minnesota"
$ws.Range("C1").Value = (Get-Date -Year 2022 -Month 3 -Day 18 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("C1").NumberFormat = "mm-dd-yy"
